$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the year label and copy formatting (style) from the row above
$ws.Range("A12").Value = "2021年"
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the numeric data for the new row (row 12)
$ws.Range("B12").Value = 36
$ws.Range("D12").Value = 39
$ws.Range("F12").Value = 2870
$ws.Range("H12").Value = 103
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = 18
$ws.Range("N12").Value = 518
$ws.Range("Q12").Value = 601
$ws.Range("R12").Value = 1
$ws.Range("T12").Value = 10
$ws.Range("W12").Value = 67
$ws.Range("Y12").Value = 1454
